$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose target values are date-like / numeric-like / percent-like
# strings ("2011-04-01", "65", "100%", ...) must be force-formatted as Text
# *before* the value is assigned, otherwise Excel auto-converts them to a
# date serial, a number, or a percentage - the target data keeps them as
# literal text. (A multi-area "A1,B2" Range only applies property changes
# to its first area in this host, so each cell is formatted individually.)

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = 'VGL_CAO_2011_2013_ Doon_tekst_200212.json'
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = '2011-04-01'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '2013-04-01'
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = '2011-10-10'
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = '4 weeks'
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = '2011-10-01'
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = 'month'
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = '2011-10-01'
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = 'o''clock'
$ws.Range("X2").NumberFormat = "@"
$ws.Range("X2").Value = '2011-10-10'
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = '4 weeks'
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = '2012-01-02'
$ws.Range("AC2").NumberFormat = "@"
$ws.Range("AC2").Value = '4 weeks'
$ws.Range("AF2").NumberFormat = "@"
$ws.Range("AF2").Value = '2012-01-01'
$ws.Range("AG2").NumberFormat = "@"
$ws.Range("AG2").Value = 'month'
$ws.Range("AJ2").NumberFormat = "@"
$ws.Range("AJ2").Value = '2012-01-02'
$ws.Range("AK2").NumberFormat = "@"
$ws.Range("AK2").Value = 'o''clock'
$ws.Range("AM2").NumberFormat = "@"
$ws.Range("AM2").Value = 'A'
$ws.Range("AN2").NumberFormat = "@"
$ws.Range("AN2").Value = 'Wages per age and function years. Employees aged 13 or 14 receive the wages from a 15-year-old. There is a separate, lower inflow scale for the long -term unemployed for a maximum of one year.'
$ws.Range("AO2").NumberFormat = "@"
$ws.Range("AO2").Value = '21 years and older'
$ws.Range("AP2").NumberFormat = "@"
$ws.Range("AP2").Value = 'As of January 1, 2012, the pension premium for the 2012 calendar year will be increased to 18.5%. In 2011 the premium was 18.3%.'
$ws.Range("AR2").NumberFormat = "@"
$ws.Range("AR2").Value = '65'
$ws.Range("AS2").NumberFormat = "@"
$ws.Range("AS2").Value = 'The intended retirement age is gradually increased to 62 years and four months.'
$ws.Range("AT2").NumberFormat = "@"
$ws.Range("AT2").Value = '21'
$ws.Range("AU2").NumberFormat = "@"
$ws.Range("AU2").Value = '2 working days'
$ws.Range("AV2").NumberFormat = "@"
$ws.Range("AV2").Value = '100%'
$ws.Range("AW2").NumberFormat = "@"
$ws.Range("AW2").Value = 'Concerns maternity leave for the partner: during delivery and two working days leave while retaining wages. Reference is made to the Work and Care Act for maternity and maternity leave.'
$ws.Range("AY2").NumberFormat = "@"
$ws.Range("AY2").Value = 'working days'
$ws.Range("AZ2").NumberFormat = "@"
$ws.Range("AZ2").Value = 'Extra vacation days based on service time (3 days at 25-40 years of service, 5 days at 40+ years of service) or age (1 day at 50-55 years, 2 days by 55-60 years, 4 days by 60+ years). The highest number of the two schemes applies.'
$ws.Range("BA2").NumberFormat = "@"
$ws.Range("BA2").Value = 'at least one month for monthly laborers, at least four weeks for periods of periods'
$ws.Range("BB2").NumberFormat = "@"
$ws.Range("BB2").Value = 'month / weeks'
$ws.Range("BC2").NumberFormat = "@"
$ws.Range("BC2").Value = 'Termination takes place at the end of a wage payment period. Permission from CWI necessary if the employee does not agree.'
$ws.Range("BD2").NumberFormat = "@"
$ws.Range("BD2").Value = 'at least one month for monthly laborers, at least four weeks for periods of periods'
$ws.Range("BE2").NumberFormat = "@"
$ws.Range("BE2").Value = 'month / weeks'
$ws.Range("BF2").NumberFormat = "@"
$ws.Range("BF2").Value = 'Termination takes place at the end of a wage payment period.'
$ws.Range("BG2").NumberFormat = "@"
$ws.Range("BG2").Value = '2'
$ws.Range("BH2").NumberFormat = "@"
$ws.Range("BH2").Value = 'months'
$ws.Range("BI2").NumberFormat = "@"
$ws.Range("BI2").Value = 'Applies to every new employment. Can be shortened or omitted in writing.'
$ws.Range("BJ2").NumberFormat = "@"
$ws.Range("BJ2").Value = 'on average 40 hours a week'
$ws.Range("BK2").NumberFormat = "@"
$ws.Range("BK2").Value = 'normal hourly wage plus a surcharge of 35%'
$ws.Range("BL2").NumberFormat = "@"
$ws.Range("BL2").Value = '9 hours a day. A maximum overtime of 3 hours a day and 10 hours a week applies to managers.'
$ws.Range("BM2").NumberFormat = "@"
$ws.Range("BM2").Value = 'Auxiliary forces are employed and paid for at least two hours a week.'
$ws.Range("BN2").NumberFormat = "@"
$ws.Range("BN2").Value = 'Shops: Mon-VR 00: 00-06: 00 (50%), 20: 00-21: 00 (33 1/3%), 21: 00-24: 00 (50%); Sat 18: 00-24: 00 (50%); Sun/holiday (100%). Distribution centers: Mon-Fri 20: 00-22: 00 (25%), 22: 00-06: 00 (50%); Sat 06: 00-18: 00 (30%), 18: 00-06: 00 (50%); Sun/holiday (100%). Freeze cell allowance DC: 8% on the hourly wage.'
$ws.Range("BQ2").NumberFormat = "@"
$ws.Range("BQ2").Value = 'Employee is informed annually about training options. Agreements about training are recorded in writing. For professional drivers (code 95), course costs, exam fees, travel costs and course time are reimbursed by the employer, with a possible repayment scheme on departure within 3 years.'

# Numeric columns
$ws.Range("K2").Value = 499.37
$ws.Range("N2").Value = 1.8
$ws.Range("O2").Value = 543.0700000000001
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 3.12
$ws.Range("V2").Value = 1.5
$ws.Range("W2").Value = 464.07
$ws.Range("Z2").Value = 1.5
$ws.Range("AA2").Value = 506.86
$ws.Range("AE2").Value = 551.21
$ws.Range("AI2").Value = 3.17
$ws.Range("AX2").Value = 24

Write-Output "applied edits"
